$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 7353799.5
$ws.Range("I33").Value = 11363853
$ws.Range("K33").Value = 11363853
$ws.Range("M33").Value = -11363624
# Row 62
$ws.Range("H62").Value = 4759.375
$ws.Range("I62").Value = 4518
$ws.Range("K62").Value = 4518
$ws.Range("M62").Value = -3894
# Row 65
$ws.Range("H65").Value = 4759.375
$ws.Range("I65").Value = 4518
$ws.Range("K65").Value = 22590
$ws.Range("M65").Value = -19470
# Row 74
$ws.Range("H74").Value = 4421.5
$ws.Range("I74").Value = 3452.5
$ws.Range("J74").Value = 5875
$ws.Range("K74").Value = 3452.5
$ws.Range("L74").Value = 5875
$ws.Range("M74").Value = -2516.5
$ws.Range("N74").Value = -7747
# Row 77
$ws.Range("H77").Value = 4421.5
$ws.Range("I77").Value = 3452.5
$ws.Range("J77").Value = 5875
$ws.Range("K77").Value = 17262.5
$ws.Range("L77").Value = 29375
$ws.Range("M77").Value = -12582.5
$ws.Range("N77").Value = -38735
# Row 116
$ws.Range("H116").Value = 4987.9
$ws.Range("I116").Value = 4987.9
$ws.Range("K116").Value = 4987.9
$ws.Range("M116").Value = -1545.9
# Row 132
$ws.Range("H132").Value = 6095.3486
$ws.Range("I132").Value = 4060.6667
$ws.Range("K132").Value = 12182.0001
$ws.Range("M132").Value = -9652.000100000001
# Row 137
$ws.Range("H137").Value = 8095.1113
$ws.Range("I137").Value = 3112.7778
$ws.Range("J137").Value = 13077.444
$ws.Range("K137").Value = 9338.3334
$ws.Range("L137").Value = 39232.33199999999
$ws.Range("M137").Value = -6788.3334
$ws.Range("N137").Value = -44332.33199999999
# Row 138
$ws.Range("H138").Value = 4110.3584
$ws.Range("I138").Value = 4335.923
$ws.Range("J138").Value = 4056.0557
$ws.Range("K138").Value = 13007.769
$ws.Range("L138").Value = 12168.1671
$ws.Range("M138").Value = -7867.769
$ws.Range("N138").Value = -22448.1671
# Row 141
$ws.Range("H141").Value = 4373.4614
$ws.Range("J141").Value = 5098.75
$ws.Range("L141").Value = 15296.25
$ws.Range("N141").Value = -25656.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 16
$ws.Range("H16").Value = 1414.25
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 1385.6666
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1385.6666
$ws.Range("M16").Value = -1213
$ws.Range("N16").Value = -1959.6666
# Row 45
$ws.Range("H45").Value = 2039.55
$ws.Range("I45").Value = 1870.3529
$ws.Range("J45").Value = 2998.3333
$ws.Range("K45").Value = 1870.3529
$ws.Range("L45").Value = 2998.3333
$ws.Range("M45").Value = -1493.3529
$ws.Range("N45").Value = -3752.3333
# Row 61
$ws.Range("H61").Value = 8954.441000000001
$ws.Range("I61").Value = 5741.174
$ws.Range("K61").Value = 5741.174
$ws.Range("M61").Value = -5529.174
# Row 102
$ws.Range("H102").Value = 1209
$ws.Range("I102").Value = 972.6
$ws.Range("K102").Value = 972.6
$ws.Range("M102").Value = 649.4
# Row 122
$ws.Range("H122").Value = 1541423.2
$ws.Range("I122").Value = 1669458.6
$ws.Range("K122").Value = 5008375.800000001
$ws.Range("M122").Value = -5005925.800000001
# Row 132
$ws.Range("H132").Value = 2580420.8
$ws.Range("I132").Value = 9889.433999999999
$ws.Range("J132").Value = 11148859
$ws.Range("K132").Value = 29668.302
$ws.Range("L132").Value = 33446577
$ws.Range("M132").Value = -27138.302
$ws.Range("N132").Value = -33451637
# Row 136
$ws.Range("H136").Value = 8954.441000000001
$ws.Range("I136").Value = 5741.174
$ws.Range("K136").Value = 17223.522
$ws.Range("M136").Value = -14673.522

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1587.0526
$ws.Range("I80").Value = 493
$ws.Range("K80").Value = 493
$ws.Range("M80").Value = 505
# Row 83
$ws.Range("H83").Value = 1587.0526
$ws.Range("I83").Value = 493
$ws.Range("K83").Value = 2465
$ws.Range("M83").Value = 2527
# Row 86
$ws.Range("H86").Value = 1332.3334
$ws.Range("J86").Value = 1165.3334
$ws.Range("L86").Value = 1165.3334
$ws.Range("N86").Value = -3411.3334
# Row 89
$ws.Range("H89").Value = 1332.3334
$ws.Range("J89").Value = 1165.3334
$ws.Range("L89").Value = 5826.666999999999
$ws.Range("N89").Value = -17058.667
# Row 105
$ws.Range("H105").Value = 1793
$ws.Range("I105").Value = 1846.1666
$ws.Range("J105").Value = 1474
$ws.Range("K105").Value = 1846.1666
$ws.Range("L105").Value = 1474
$ws.Range("M105").Value = -99.16660000000002
$ws.Range("N105").Value = -4968
# Row 134
$ws.Range("H134").Value = 16430.482
$ws.Range("I134").Value = 10373.3
$ws.Range("K134").Value = 31119.9
$ws.Range("M134").Value = -28584.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 610.41174
$ws.Range("J7").Value = 1551.3334
$ws.Range("L7").Value = 1551.3334
$ws.Range("N7").Value = -1777.3334
# Row 16
$ws.Range("H16").Value = 10550.643
$ws.Range("J16").Value = 15886.125
$ws.Range("L16").Value = 15886.125
$ws.Range("N16").Value = -16460.125
# Row 22
$ws.Range("H22").Value = 2249.1667
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 2599
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 2599
$ws.Range("M22").Value = -150
$ws.Range("N22").Value = -3299
# Row 105
$ws.Range("H105").Value = 9764.083000000001
$ws.Range("I105").Value = 15310
$ws.Range("K105").Value = 15310
$ws.Range("M105").Value = -13563
# Row 113
$ws.Range("H113").Value = 10550.643
$ws.Range("J113").Value = 15886.125
$ws.Range("L113").Value = 15886.125
$ws.Range("N113").Value = -20226.125
# Row 122
$ws.Range("H122").Value = 1543.7
$ws.Range("I122").Value = 1105
$ws.Range("K122").Value = 3315
$ws.Range("M122").Value = -865

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 7772098.5
$ws.Range("I32").Value = 253036.75
$ws.Range("J32").Value = 11113904
$ws.Range("K32").Value = 759110.25
$ws.Range("L32").Value = 33341712
$ws.Range("M32").Value = -758827.25
$ws.Range("N32").Value = -33342278
# Row 113
$ws.Range("H113").Value = 6060.227
$ws.Range("I113").Value = 15363.625
$ws.Range("K113").Value = 46090.875
$ws.Range("M113").Value = -43920.875
# Row 122
$ws.Range("H122").Value = 15372828
$ws.Range("J122").Value = 4047761.2
$ws.Range("L122").Value = 36429850.8
$ws.Range("N122").Value = -36434750.8
# Row 129
$ws.Range("H129").Value = 8138.7812
$ws.Range("J129").Value = 3416.2
$ws.Range("L129").Value = 10248.6
$ws.Range("N129").Value = -20248.6
# Row 131
$ws.Range("H131").Value = 1487.92
$ws.Range("J131").Value = 1487.92
$ws.Range("L131").Value = 4463.76
$ws.Range("N131").Value = -14543.76

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 6444.643
$ws.Range("I102").Value = 6825.077
$ws.Range("K102").Value = 6825.077
$ws.Range("M102").Value = -5203.077
# Row 113
$ws.Range("H113").Value = 2700.7273
$ws.Range("I113").Value = 2670.8
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2670.8
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -500.8000000000002
$ws.Range("N113").Value = -7340
# Row 132
$ws.Range("H132").Value = 5799.6665
$ws.Range("I132").Value = 4191.5
$ws.Range("K132").Value = 12574.5
$ws.Range("M132").Value = -10044.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6536.6816
$ws.Range("I40").Value = 6050.4375
$ws.Range("K40").Value = 6050.4375
$ws.Range("M40").Value = -5914.4375
# Row 61
$ws.Range("H61").Value = 6107.154
$ws.Range("I61").Value = 4940.6
$ws.Range("K61").Value = 4940.6
$ws.Range("M61").Value = -4738.6
# Row 93
$ws.Range("H93").Value = 7398.75
$ws.Range("I93").Value = 7560.737
$ws.Range("K93").Value = 7560.737
$ws.Range("M93").Value = -6312.737
# Row 100
$ws.Range("H100").Value = 4678.391
$ws.Range("J100").Value = 3846
$ws.Range("L100").Value = 3846
$ws.Range("N100").Value = -4928
# Row 113
$ws.Range("H113").Value = 6107.154
$ws.Range("I113").Value = 4940.6
$ws.Range("K113").Value = 4940.6
$ws.Range("M113").Value = -2770.6
# Row 122
$ws.Range("H122").Value = 7755.407
$ws.Range("I122").Value = 7249.5
$ws.Range("J122").Value = 8767.223
$ws.Range("K122").Value = 21748.5
$ws.Range("L122").Value = 26301.669
$ws.Range("M122").Value = -19298.5
$ws.Range("N122").Value = -31201.669
# Row 132
$ws.Range("H132").Value = 898705.9399999999
$ws.Range("I132").Value = 4345.3
$ws.Range("J132").Value = 1614194.5
$ws.Range("K132").Value = 13035.9
$ws.Range("L132").Value = 4842583.5
$ws.Range("M132").Value = -10505.9
$ws.Range("N132").Value = -4847643.5
# Row 139
$ws.Range("H139").Value = 42500
$ws.Range("J139").Value = 42500
$ws.Range("L139").Value = 42500
$ws.Range("N139").Value = -52780

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 3130.111
$ws.Range("J23").Value = 6794.25
$ws.Range("L23").Value = 6794.25
$ws.Range("N23").Value = -7252.25
# Row 100
$ws.Range("H100").Value = 446.53845
$ws.Range("I100").Value = 423.81818
$ws.Range("J100").Value = 571.5
$ws.Range("K100").Value = 847.63636
$ws.Range("L100").Value = 1143
$ws.Range("M100").Value = -306.63636
$ws.Range("N100").Value = -2225
